$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert 4 new rows at the top; this shifts the existing 72 rows down to 5-76
#    and also duplicates the old H1:H8 formula cells down into H5:H12.
$ws.Rows("1:4").Insert()

# 2. Those duplicated H5:H12 cells should no longer carry the generated-array
#    formula (only the 4 newest rows keep it), so clear them out entirely.
$ws.Range("H5:H12").ClearContents()

# 3. Give the 4 new rows the same number formats as the existing data rows
#    (date in col A, text in col D, text in cols F/G) by copying formats down
#    from row 5 (the row that used to be row 1).
$ws.Range("A5").Copy()
$ws.Range("A1:A4").PasteSpecial(-4122)
$ws.Range("D5").Copy()
$ws.Range("D1:D4").PasteSpecial(-4122)
$ws.Range("F5:G5").Copy()
$ws.Range("F1:G4").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# 4. Populate the 4 new transaction rows (statement entries for 2014-03-05).
$ws.Range("A1").Value = 41703
$ws.Range("B1").Value = "  TRANSFERENCIA INTERNET"
$ws.Range("C1").Value = "D"
$ws.Range("D1").Value = "0016349162"
$ws.Range("E1").Value = "AG. NORTE"
$ws.Range("F1").Value = "280.95  "
$ws.Range("G1").Value = "187.30"

$ws.Range("A2").Value = 41703
$ws.Range("B2").Value = "CONSUMO VISA NA PIZZA HUT"
$ws.Range("C2").Value = "D"
$ws.Range("D2").Value = "0014037864"
$ws.Range("E2").Value = "INSTITUCIONAL SS.CC."
$ws.Range("F2").Value = "29.19  "
$ws.Range("G2").Value = "468.25"

$ws.Range("A3").Value = 41703
$ws.Range("B3").Value = "CONSUMO DATA AKI MOLINEROS 161"
$ws.Range("C3").Value = "D"
$ws.Range("D3").Value = "0011137163"
$ws.Range("E3").Value = "INSTITUCIONAL SS.CC."
$ws.Range("F3").Value = "19.58  "
$ws.Range("G3").Value = "497.44"

$ws.Range("A4").Value = 41703
$ws.Range("B4").Value = "RETIRO ATM BP N/S.S CCI-2"
$ws.Range("C4").Value = "D"
$ws.Range("D4").Value = "0007326973"
$ws.Range("E4").Value = "C.C.I"
$ws.Range("F4").Value = "50.00  "
$ws.Range("G4").Value = "517.02"

# 5. Rebuild the PHP-array-building formula for the 4 newest rows only (H1:H4).
$formula = '=CONCATENATE("array(''mo_fecha'' => new \DateTime(''",TEXT(A1,"yyyy-mm-dd"),"''), ''mo_concepto'' => ''",B1,"'', ''mo_tipo'' => ''",C1,"'', ''mo_documento'' => ''",D1,"'', ''mo_oficina'' => ''",E1,"'', ''mo_monto'' => ",F1,", ''mo_saldo'' => ",G1,", ''mo_fecha_crea'' => new \DateTime(''",TEXT(NOW(),"yyyy-mm-dd H:m:s"),"''), ''mo_quien_crea'' => 1, ''mo_fecha_modifica'' => NULL, ''mo_quien_modifica'' => NULL, ''mo_borrado_logico'' => false),")'
$ws.Range("H1:H4").Formula = $formula

# 6. Update the selection to match the new formula range.
$ws.Range("H1:H4").Select()

Write-Host "done"
